# Update Hangboard data through april 21st
# Appends five new rows (20-24) of hangboard training data to Sheet1,
# covering 25 Mar 2020 through 21 Apr 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("25 Mar 2020", "0,6", "-12.5,6", "-25,6", "-50,6", "-32.5,5,9", "-20,6", "-42.5,5,7", "-25,3,9,7,6", "-27.5,4,7,8"),
    @("30 Mar 2020", "0,6", "-12.5,6", "-25,6", "-50,6", "-32.5,6", "-20,6", "-42.5,6", "-25,6", "-27.5,6"),
    @("9 Apr 2020", "0,6", "-10,6", "-22.5,6", "-47.5,6", "-27.5,6", "-17.5,6", "-40,6", "-22.5,5,7", "-25,5,8"),
    @("15 Apr 2020", "0,6", "-10,6", "-22.5,6", "-47.5,6", "-27.5,5,9", "-17.5,4,9,5", "-40,4,8,7", "-22.5,3,5,5,5", "-25,4,8,8"),
    @("21 Apr 2020", "0,6", "-10,6", "-22.5,6", "-47.5,6", "-27.5,5,9", "-17.5,4,9,7", "-40,5,8", "-22.5,3,6,7,5", "-25,3,7,7,7")
)

$startRow = 20
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

# Match the new active selection left by the editing session.
$ws.Range("J27").Select()
